# "fixed some errors in the result files"
# The tournament schedule was compressed (group-phase & knockout matches
# start/end earlier than originally listed) and two typos in a player's
# name were corrected ("Bertalan Pusztai" -> "Berci Pusztai").
# Only the affected Starting time / End time / Team cells are touched;
# everything else (results, tables, groups) is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Name correction: "Bertalan Pusztai" -> "Berci Pusztai" ---
$ws.Range("F2").Value  = "Berci Pusztai"
$ws.Range("F6").Value  = "Berci Pusztai"
$ws.Range("F10").Value = "Berci Pusztai"
$ws.Range("F14").Value = "Berci Pusztai"
$ws.Range("F18").Value = "Berci Pusztai"
$ws.Range("F22").Value = "Berci Pusztai"
$ws.Range("F26").Value = "Berci Pusztai"
$ws.Range("F31").Value = "Berci Pusztai"
$ws.Range("H32").Value = "Berci Pusztai"

# --- Schedule fix: shift starting / end times for rows 2-33 ---
# Row 2-3 (12:00 start is unchanged, only end time moves)
$ws.Range("B2").Value  = "12:25"
$ws.Range("B3").Value  = "12:25"

# Row 4-5
$ws.Range("A4").Value  = "12:30"
$ws.Range("B4").Value  = "12:55"
$ws.Range("A5").Value  = "12:30"
$ws.Range("B5").Value  = "12:55"

# Row 6-7
$ws.Range("A6").Value  = "13:00"
$ws.Range("B6").Value  = "13:25"
$ws.Range("A7").Value  = "13:00"
$ws.Range("B7").Value  = "13:25"

# Row 8-9
$ws.Range("A8").Value  = "13:30"
$ws.Range("B8").Value  = "13:55"
$ws.Range("A9").Value  = "13:30"
$ws.Range("B9").Value  = "13:55"

# Row 10-11
$ws.Range("A10").Value = "14:00"
$ws.Range("B10").Value = "14:25"
$ws.Range("A11").Value = "14:00"
$ws.Range("B11").Value = "14:25"

# Row 12-13
$ws.Range("A12").Value = "14:30"
$ws.Range("B12").Value = "14:55"
$ws.Range("A13").Value = "14:30"
$ws.Range("B13").Value = "14:55"

# Row 14-15
$ws.Range("A14").Value = "15:00"
$ws.Range("B14").Value = "15:25"
$ws.Range("A15").Value = "15:00"
$ws.Range("B15").Value = "15:25"

# Row 16-17
$ws.Range("A16").Value = "15:30"
$ws.Range("B16").Value = "15:55"
$ws.Range("A17").Value = "15:30"
$ws.Range("B17").Value = "15:55"

# Row 18-19
$ws.Range("A18").Value = "16:00"
$ws.Range("B18").Value = "16:25"
$ws.Range("A19").Value = "16:00"
$ws.Range("B19").Value = "16:25"

# Row 20-21
$ws.Range("A20").Value = "16:30"
$ws.Range("B20").Value = "16:55"
$ws.Range("A21").Value = "16:30"
$ws.Range("B21").Value = "16:55"

# Row 22-23
$ws.Range("A22").Value = "17:00"
$ws.Range("B22").Value = "17:25"
$ws.Range("A23").Value = "17:00"
$ws.Range("B23").Value = "17:25"

# Row 24-25
$ws.Range("A24").Value = "17:30"
$ws.Range("B24").Value = "17:55"
$ws.Range("A25").Value = "17:30"
$ws.Range("B25").Value = "17:55"

# Row 26-27
$ws.Range("A26").Value = "18:00"
$ws.Range("B26").Value = "18:25"
$ws.Range("A27").Value = "18:00"
$ws.Range("B27").Value = "18:25"

# Row 28-29
$ws.Range("A28").Value = "18:30"
$ws.Range("B28").Value = "18:55"
$ws.Range("A29").Value = "18:30"
$ws.Range("B29").Value = "18:55"

# Row 30-31 (Knockout phase - Match B1 / Match B2)
$ws.Range("A30").Value = "19:00"
$ws.Range("B30").Value = "19:25"
$ws.Range("A31").Value = "19:00"
$ws.Range("B31").Value = "19:25"

# Row 32-33 (Knockout phase - Match B3 / Match B4)
$ws.Range("A32").Value = "19:30"
$ws.Range("B32").Value = "19:55"
$ws.Range("A33").Value = "19:30"
$ws.Range("B33").Value = "19:55"
